# Scheduled-runner data refresh for the Ragnarok_Profits workbook.
# Market-price driven columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# are rewritten per leve/row with freshly pulled averages; everything else
# (names, ids, levels, static gil/exp) is left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 1265.2727
$ws.Range("J4").Value = 1899.7142
$ws.Range("L4").Value = 1899.7142
$ws.Range("N4").Value = -2127.7142

# Row 5: Met a Sticky End
$ws.Range("H5").Value = 23.333334
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 17: One for the Road
$ws.Range("H17").Value = 26393.871
$ws.Range("J17").Value = 28560.027
$ws.Range("L17").Value = 85680.08099999999
$ws.Range("N17").Value = -86016.08099999999

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 2458.4
$ws.Range("I43").Value = 3112
$ws.Range("J43").Value = 2022.6666
$ws.Range("K43").Value = 3112
$ws.Range("L43").Value = 2022.6666
$ws.Range("M43").Value = -3043
$ws.Range("N43").Value = -2160.6666

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 1861.375
$ws.Range("I107").Value = 611.7
$ws.Range("K107").Value = 611.7
$ws.Range("M107").Value = 1308.3

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 5990.5312
$ws.Range("I132").Value = 2408.6667
$ws.Range("K132").Value = 7226.000100000001
$ws.Range("M132").Value = -4696.000100000001

# Row 135: For Tired Minds
$ws.Range("H135").Value = 1598.5667
$ws.Range("J135").Value = 8249
$ws.Range("L135").Value = 74241
$ws.Range("N135").Value = -79311

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 500239.2
$ws.Range("I137").Value = 646
$ws.Range("J137").Value = 1443915.2
$ws.Range("K137").Value = 1938
$ws.Range("L137").Value = 4331745.6
$ws.Range("M137").Value = 612
$ws.Range("N137").Value = -4336845.6

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2791.869
$ws.Range("J138").Value = 4916.304
$ws.Range("L138").Value = 14748.912
$ws.Range("N138").Value = -25028.912

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4298.7637
$ws.Range("I32").Value = 4315.245
$ws.Range("K32").Value = 4315.245
$ws.Range("M32").Value = -4028.245

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2385041
$ws.Range("I132").Value = 4005.9211
$ws.Range("K132").Value = 12017.7633
$ws.Range("M132").Value = -9487.763300000001

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 483560.9
$ws.Range("I86").Value = 848827
$ws.Range("K86").Value = 848827
$ws.Range("M86").Value = -847704

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 483560.9
$ws.Range("I89").Value = 848827
$ws.Range("K89").Value = 4244135
$ws.Range("M89").Value = -4238519

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 3437.0908
$ws.Range("I99").Value = 2738.625
$ws.Range("K99").Value = 2738.625
$ws.Range("M99").Value = -1240.625

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1210.75
$ws.Range("J22").Value = 1617.4
$ws.Range("L22").Value = 1617.4
$ws.Range("N22").Value = -2317.4

# Row 31: Wall Not Found
$ws.Range("H31").Value = 25902468
$ws.Range("I31").Value = 40002404
$ws.Range("K31").Value = 40002404
$ws.Range("M31").Value = -40002109

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 25902468
$ws.Range("I34").Value = 40002404
$ws.Range("K34").Value = 40002404
$ws.Range("M34").Value = -40002202

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2662.75
$ws.Range("I58").Value = 2632.647
$ws.Range("J58").Value = 2833.3333
$ws.Range("K58").Value = 2632.647
$ws.Range("L58").Value = 2833.3333
$ws.Range("M58").Value = -2429.647
$ws.Range("N58").Value = -3239.3333

# Row 99: O Pine
$ws.Range("H99").Value = 14714.889
$ws.Range("I99").Value = 7701.231
$ws.Range("K99").Value = 7701.231
$ws.Range("M99").Value = -6203.231

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 4194.722
$ws.Range("I122").Value = 4362.7
$ws.Range("J122").Value = 3984.75
$ws.Range("K122").Value = 13088.1
$ws.Range("L122").Value = 11954.25
$ws.Range("M122").Value = -10638.1
$ws.Range("N122").Value = -16854.25

# Row 126: A Better Conductor
$ws.Range("H126").Value = 14714.889
$ws.Range("I126").Value = 7701.231
$ws.Range("K126").Value = 23103.693
$ws.Range("M126").Value = -20633.693

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2554.4443
$ws.Range("I132").Value = 2316.1304
$ws.Range("K132").Value = 6948.3912
$ws.Range("M132").Value = -4418.3912

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3461.05
$ws.Range("I134").Value = 3327.125
$ws.Range("J134").Value = 3996.75
$ws.Range("K134").Value = 9981.375
$ws.Range("L134").Value = 11990.25
$ws.Range("M134").Value = -7446.375
$ws.Range("N134").Value = -17060.25

# Row 136: Turali Quality
$ws.Range("H136").Value = 2662.75
$ws.Range("I136").Value = 2632.647
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 7897.941
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -5347.941
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On
$ws.Range("H11").Value = 83336750
$ws.Range("I11").Value = 525.8889
$ws.Range("K11").Value = 1577.6667
$ws.Range("M11").Value = -1437.6667

# Row 33: Cooking with Gas
$ws.Range("H33").Value = 6299864.5
$ws.Range("I33").Value = 193.75
$ws.Range("J33").Value = 9899676
$ws.Range("K33").Value = 1162.5
$ws.Range("L33").Value = 59398056
$ws.Range("M33").Value = -879.5
$ws.Range("N33").Value = -59398622

# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 7770.9414
$ws.Range("I134").Value = 2769.5
$ws.Range("K134").Value = 8308.5
$ws.Range("M134").Value = -3238.5

# Row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 11205.529
$ws.Range("I138").Value = 8751.571
$ws.Range("K138").Value = 26254.713
$ws.Range("M138").Value = -21114.713

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 23.4
$ws.Range("I2").Value = 32
$ws.Range("K2").Value = 32
$ws.Range("M2").Value = 81

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2624.1853
$ws.Range("J80").Value = 4481.9
$ws.Range("L80").Value = 4481.9
$ws.Range("N80").Value = -6477.9

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2624.1853
$ws.Range("J83").Value = 4481.9
$ws.Range("L83").Value = 22409.5
$ws.Range("N83").Value = -32393.5

# Row 139: Ringing Gratitude
$ws.Range("H139").Value = 112520.22
$ws.Range("J139").Value = 112520.22
$ws.Range("L139").Value = 112520.22
$ws.Range("N139").Value = -122800.22

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 4891.3887
$ws.Range("I16").Value = 1795.4166
$ws.Range("J16").Value = 11083.333
$ws.Range("K16").Value = 1795.4166
$ws.Range("L16").Value = 11083.333
$ws.Range("M16").Value = -1625.4166
$ws.Range("N16").Value = -11423.333

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 36741
$ws.Range("I22").Value = 55056
$ws.Range("J22").Value = 111
$ws.Range("K22").Value = 55056
$ws.Range("L22").Value = 111
$ws.Range("M22").Value = -54761
$ws.Range("N22").Value = -701

# Row 27: Fire and Hide
$ws.Range("H27").Value = 36741
$ws.Range("I27").Value = 55056
$ws.Range("J27").Value = 111
$ws.Range("K27").Value = 55056
$ws.Range("L27").Value = 111
$ws.Range("M27").Value = -54949
$ws.Range("N27").Value = -325

# Row 40: Best Served Toad
$ws.Range("H40").Value = 4436.3
$ws.Range("I40").Value = 4349.875
$ws.Range("J40").Value = 4782
$ws.Range("K40").Value = 4349.875
$ws.Range("L40").Value = 4782
$ws.Range("M40").Value = -4213.875
$ws.Range("N40").Value = -5054

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1493.5
$ws.Range("I46").Value = 1324.6666
$ws.Range("K46").Value = 1324.6666
$ws.Range("M46").Value = -1136.6666

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 1323.2
$ws.Range("I55").Value = 637
$ws.Range("K55").Value = 637
$ws.Range("M55").Value = -464

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 3790364.2
$ws.Range("I68").Value = 6946472.5
$ws.Range("J68").Value = 3034.6
$ws.Range("K68").Value = 6946472.5
$ws.Range("L68").Value = 3034.6
$ws.Range("M68").Value = -6945723.5
$ws.Range("N68").Value = -4532.6

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 3790364.2
$ws.Range("I71").Value = 6946472.5
$ws.Range("J71").Value = 3034.6
$ws.Range("K71").Value = 34732362.5
$ws.Range("L71").Value = 15173
$ws.Range("M71").Value = -34728618.5
$ws.Range("N71").Value = -22661

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 16997.5
$ws.Range("I96").Value = 14330
$ws.Range("K96").Value = 14330
$ws.Range("M96").Value = -12957

# Row 107: Flax Wax
$ws.Range("H107").Value = 3022.681
$ws.Range("I107").Value = 1303.3043
$ws.Range("K107").Value = 3909.9129
$ws.Range("M107").Value = -1989.9129

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 5087.1577
$ws.Range("I122").Value = 4989.364
$ws.Range("K122").Value = 14968.092
$ws.Range("M122").Value = -12518.092

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 456257.1
$ws.Range("I132").Value = 1803.4706
$ws.Range("K132").Value = 5410.4118
$ws.Range("M132").Value = -2880.4118

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 288294.16
$ws.Range("I136").Value = 2675.6365
$ws.Range("K136").Value = 8026.9095
$ws.Range("M136").Value = -5476.9095

# Row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 83128.5
$ws.Range("J138").Value = 98885.5
$ws.Range("L138").Value = 98885.5
$ws.Range("N138").Value = -109165.5
